# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.817.53'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.507.24'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''581.15'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').Value = '''161.99'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '''0.602'
$ws.Range('E8').Value = '  +10.16%  '
$ws.Range('D9').Value = '3.511.57'
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('D10').Value = '''7.32'
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  +2.44%  '
$ws.Range('D12').Value = '''0.447'
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('D13').Value = '4.109.02'
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('E16').Value = '  +5.54%  '
$ws.Range('D17').Value = '65.817.83'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = '3.520.03'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').Value = '''6.48'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').Value = '''14.35'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('D21').Value = '''391.22'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = '''8.31'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').Value = '''73.65'
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').Value = '''0.997'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').Value = '''0.0000126'
$ws.Range('E26').Value = '  +5.23%  '
$ws.Range('D27').Value = '''9.97'
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('D28').Value = '''0.180'
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +6.52%  '
$ws.Range('E31').Value = '  +6.94%  '
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('D33').Value = '''23.80'
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').Value = '''6.51'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''7.19'
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''1.57'
$ws.Range('E36').Value = '  +7.72%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''163.10'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '''1.95'
$ws.Range('E38').Value = '  +6.56%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.094.43'
$ws.Range('E39').Value = '  +6.03%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.0776'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''27.54'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0324'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''4.55'
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''43.23'
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '''0.785'
$ws.Range('E45').Value = '  +1.75%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '''25.85'
$ws.Range('E46').Value = '  +9.47%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = '''1.12'
$ws.Range('E47').Value = '  +4.54%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '''2.27'
$ws.Range('E48').Value = '  +4.66%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '''315.77'
$ws.Range('E49').Value = '  +8.69%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''6.73'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.108'
$ws.Range('E51').Value = '  +4.62%  '
